# Update the "Förändrad" (Changed) date column C for all existing data rows
# (rows 2 through 266) from 2023-09-23 (45192) to 2023-10-03 (45202).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C266").Value = 45202

# Row 266 gains an explicit row height (15, customHeight) in the updated file.
$ws.Rows.Item(266).RowHeight = 15

# Append three new records (rows 267-269) to the bottom of the sheet.
$newRows = @(
    @{ Row = 267; A = "A 46728-2023"; B = 45198; C = 45202; G = 1.8 },
    @{ Row = 268; A = "A 46835-2023"; B = 45199; C = 45202; G = 0 },
    @{ Row = 269; A = "A 46836-2023"; B = 45199; C = 45202; G = 0.1 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A

    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($row, 4).Value = "SKÅNE LÄN"
    $ws.Cells.Item($row, 5).Value = "KLIPPAN"

    # Column F (Markägare) stays empty, as in every other row of this sheet.

    $ws.Cells.Item($row, 7).Value = $r.G

    # H (8) through Q (17) are all zero for these new records.
    $ws.Range($ws.Cells.Item($row, 8), $ws.Cells.Item($row, 17)).Value = 0

    # Column R (Artnamn) is present but empty, with the wrap-text style used
    # throughout the sheet.
    $ws.Cells.Item($row, 18).WrapText = $true
}

# Rows 267 and 268 also carry the explicit 15pt custom row height; row 269
# (the new last row) does not.
$ws.Rows.Item(267).RowHeight = 15
$ws.Rows.Item(268).RowHeight = 15

Write-Output "done"
